$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $fVal = $ws.Cells.Item($r, 6).Value2

    # Rotate columns D (group), E (category), F (category-code)
    # so that: newD = oldF, newE = oldD, newF = oldE
    $ws.Cells.Item($r, 4).Value2 = $fVal
    $ws.Cells.Item($r, 5).Value2 = $dVal
    $ws.Cells.Item($r, 6).Value2 = $eVal
}
